$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 5611.1665
$ws.Cells.Item(40, 9).Value = 8574
$ws.Cells.Item(40, 10).Value = 2648.3333
$ws.Cells.Item(40, 11).Value = 8574
$ws.Cells.Item(40, 12).Value = 2648.3333
$ws.Cells.Item(40, 13).Value = -8399
$ws.Cells.Item(40, 14).Value = -2998.3333

$ws.Cells.Item(45, 8).Value = 1250
$ws.Cells.Item(45, 9).Value = 1000
$ws.Cells.Item(45, 10).Value = 1500
$ws.Cells.Item(45, 11).Value = 3000
$ws.Cells.Item(45, 12).Value = 4500
$ws.Cells.Item(45, 13).Value = -2808
$ws.Cells.Item(45, 14).Value = -4884

$ws.Cells.Item(64, 8).Value = 8434.571
$ws.Cells.Item(64, 9).Value = 9287.5
$ws.Cells.Item(64, 10).Value = 7297.3335
$ws.Cells.Item(64, 11).Value = 9287.5
$ws.Cells.Item(64, 12).Value = 7297.3335
$ws.Cells.Item(64, 13).Value = -9039.5
$ws.Cells.Item(64, 14).Value = -7793.3335

$ws.Cells.Item(67, 8).Value = 8434.571
$ws.Cells.Item(67, 9).Value = 9287.5
$ws.Cells.Item(67, 10).Value = 7297.3335
$ws.Cells.Item(67, 11).Value = 9287.5
$ws.Cells.Item(67, 12).Value = 7297.3335
$ws.Cells.Item(67, 13).Value = -8429.5
$ws.Cells.Item(67, 14).Value = -9013.333500000001

$ws.Cells.Item(92, 8).Value = 1559.7142
$ws.Cells.Item(92, 9).Value = 1319.75
$ws.Cells.Item(92, 10).Value = 2999.5
$ws.Cells.Item(92, 11).Value = 1319.75
$ws.Cells.Item(92, 12).Value = 2999.5
$ws.Cells.Item(92, 13).Value = -71.75
$ws.Cells.Item(92, 14).Value = -5495.5

$ws.Cells.Item(97, 8).Value = 979.8
$ws.Cells.Item(97, 10).Value = 1124.75
$ws.Cells.Item(97, 12).Value = 3374.25
$ws.Cells.Item(97, 14).Value = -4366.25

$ws.Cells.Item(99, 8).Value = 953.8333
$ws.Cells.Item(99, 10).Value = 1064.6
$ws.Cells.Item(99, 12).Value = 3193.8
$ws.Cells.Item(99, 14).Value = -6189.799999999999

$ws.Cells.Item(107, 8).Value = 717.1667
$ws.Cells.Item(107, 9).Value = 717.1667
$ws.Cells.Item(107, 11).Value = 717.1667
$ws.Cells.Item(107, 13).Value = 1202.8333

$ws.Cells.Item(135, 8).Value = 4000
$ws.Cells.Item(135, 9).Value = 4000
$ws.Cells.Item(135, 11).Value = 36000
$ws.Cells.Item(135, 13).Value = -33465

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 493.8125
$ws.Cells.Item(5, 9).Value = 484.69232
$ws.Cells.Item(5, 10).Value = 533.3333
$ws.Cells.Item(5, 11).Value = 484.69232
$ws.Cells.Item(5, 12).Value = 533.3333
$ws.Cells.Item(5, 13).Value = -372.69232
$ws.Cells.Item(5, 14).Value = -757.3333

$ws.Cells.Item(45, 8).Value = 2646.6
$ws.Cells.Item(45, 9).Value = 3138.1428
$ws.Cells.Item(45, 11).Value = 3138.1428
$ws.Cells.Item(45, 13).Value = -2761.1428

$ws.Cells.Item(61, 8).Value = 62501884
$ws.Cells.Item(61, 9).Value = 83335160
$ws.Cells.Item(61, 10).Value = 2049.5
$ws.Cells.Item(61, 11).Value = 83335160
$ws.Cells.Item(61, 12).Value = 2049.5
$ws.Cells.Item(61, 13).Value = -83334948
$ws.Cells.Item(61, 14).Value = -2473.5

$ws.Cells.Item(97, 8).Value = 1471.6818
$ws.Cells.Item(97, 9).Value = 1655.3889
$ws.Cells.Item(97, 10).Value = 645
$ws.Cells.Item(97, 11).Value = 1655.3889
$ws.Cells.Item(97, 12).Value = 645
$ws.Cells.Item(97, 13).Value = -1159.3889
$ws.Cells.Item(97, 14).Value = -1637

$ws.Cells.Item(109, 8).Value = 62688.5
$ws.Cells.Item(109, 10).Value = 62688.5
$ws.Cells.Item(109, 12).Value = 62688.5
$ws.Cells.Item(109, 14).Value = -65462.5

$ws.Cells.Item(122, 8).Value = 12349895
$ws.Cells.Item(122, 9).Value = 2352.8333
$ws.Cells.Item(122, 11).Value = 7058.499899999999
$ws.Cells.Item(122, 13).Value = -4608.499899999999

$ws.Cells.Item(136, 8).Value = 62501884
$ws.Cells.Item(136, 9).Value = 83335160
$ws.Cells.Item(136, 10).Value = 2049.5
$ws.Cells.Item(136, 11).Value = 250005480
$ws.Cells.Item(136, 12).Value = 6148.5
$ws.Cells.Item(136, 13).Value = -250002930
$ws.Cells.Item(136, 14).Value = -11248.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 493.8125
$ws.Cells.Item(4, 9).Value = 484.69232
$ws.Cells.Item(4, 10).Value = 533.3333
$ws.Cells.Item(4, 11).Value = 484.69232
$ws.Cells.Item(4, 12).Value = 533.3333
$ws.Cells.Item(4, 13).Value = -369.69232
$ws.Cells.Item(4, 14).Value = -763.3333

$ws.Cells.Item(64, 8).Value = 1869.6
$ws.Cells.Item(64, 9).Value = 1782.3334
$ws.Cells.Item(64, 10).Value = 2000.5
$ws.Cells.Item(64, 11).Value = 1782.3334
$ws.Cells.Item(64, 12).Value = 2000.5
$ws.Cells.Item(64, 13).Value = -1557.3334
$ws.Cells.Item(64, 14).Value = -2450.5

$ws.Cells.Item(67, 8).Value = 1869.6
$ws.Cells.Item(67, 9).Value = 1782.3334
$ws.Cells.Item(67, 10).Value = 2000.5
$ws.Cells.Item(67, 11).Value = 1782.3334
$ws.Cells.Item(67, 12).Value = 2000.5
$ws.Cells.Item(67, 13).Value = -1002.3334
$ws.Cells.Item(67, 14).Value = -3560.5

$ws.Cells.Item(94, 8).Value = 535.0625
$ws.Cells.Item(94, 9).Value = 549.06665
$ws.Cells.Item(94, 11).Value = 549.06665
$ws.Cells.Item(94, 13).Value = -98.06664999999998

$ws.Cells.Item(135, 8).Value = 64969.6
$ws.Cells.Item(135, 10).Value = 64969.6
$ws.Cells.Item(135, 12).Value = 64969.6
$ws.Cells.Item(135, 14).Value = -75109.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 833.875
$ws.Cells.Item(69, 10).Value = 1199.5
$ws.Cells.Item(69, 12).Value = 3598.5
$ws.Cells.Item(69, 14).Value = -5220.5

$ws.Cells.Item(72, 8).Value = 833.875
$ws.Cells.Item(72, 10).Value = 1199.5
$ws.Cells.Item(72, 12).Value = 10795.5
$ws.Cells.Item(72, 14).Value = -18907.5

$ws.Cells.Item(107, 8).Value = 1202.619
$ws.Cells.Item(107, 9).Value = 907.25
$ws.Cells.Item(107, 10).Value = 2147.8
$ws.Cells.Item(107, 11).Value = 2721.75
$ws.Cells.Item(107, 12).Value = 6443.400000000001
$ws.Cells.Item(107, 13).Value = -801.75
$ws.Cells.Item(107, 14).Value = -10283.4

$ws.Cells.Item(122, 8).Value = 1454.8182
$ws.Cells.Item(122, 9).Value = 819.3333
$ws.Cells.Item(122, 10).Value = 1693.125
$ws.Cells.Item(122, 11).Value = 7373.9997
$ws.Cells.Item(122, 12).Value = 15238.125
$ws.Cells.Item(122, 13).Value = -4923.9997
$ws.Cells.Item(122, 14).Value = -20138.125

$ws.Cells.Item(131, 8).Value = 3325.923
$ws.Cells.Item(131, 10).Value = 4523.375
$ws.Cells.Item(131, 12).Value = 13570.125
$ws.Cells.Item(131, 14).Value = -23650.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 646.9167
$ws.Cells.Item(97, 9).Value = 483.3
$ws.Cells.Item(97, 11).Value = 483.3
$ws.Cells.Item(97, 13).Value = 12.69999999999999

$ws.Cells.Item(122, 8).Value = 38463988
$ws.Cells.Item(122, 9).Value = 1894.5555
$ws.Cells.Item(122, 11).Value = 5683.666499999999
$ws.Cells.Item(122, 13).Value = -3233.666499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1315.8636
$ws.Cells.Item(16, 9).Value = 1340.4286
$ws.Cells.Item(16, 11).Value = 1340.4286
$ws.Cells.Item(16, 13).Value = -1170.4286

$ws.Cells.Item(22, 8).Value = 3320
$ws.Cells.Item(22, 9).Value = 3266.6667
$ws.Cells.Item(22, 10).Value = 3400
$ws.Cells.Item(22, 11).Value = 3266.6667
$ws.Cells.Item(22, 12).Value = 3400
$ws.Cells.Item(22, 13).Value = -2971.6667
$ws.Cells.Item(22, 14).Value = -3990

$ws.Cells.Item(27, 8).Value = 3320
$ws.Cells.Item(27, 9).Value = 3266.6667
$ws.Cells.Item(27, 10).Value = 3400
$ws.Cells.Item(27, 11).Value = 3266.6667
$ws.Cells.Item(27, 12).Value = 3400
$ws.Cells.Item(27, 13).Value = -3159.6667
$ws.Cells.Item(27, 14).Value = -3614

$ws.Cells.Item(46, 8).Value = 1648.4062
$ws.Cells.Item(46, 9).Value = 989.96
$ws.Cells.Item(46, 11).Value = 989.96
$ws.Cells.Item(46, 13).Value = -801.96

$ws.Cells.Item(68, 8).Value = 5666.6665
$ws.Cells.Item(68, 10).Value = 5500
$ws.Cells.Item(68, 12).Value = 5500
$ws.Cells.Item(68, 14).Value = -6998

$ws.Cells.Item(71, 8).Value = 5666.6665
$ws.Cells.Item(71, 10).Value = 5500
$ws.Cells.Item(71, 12).Value = 27500
$ws.Cells.Item(71, 14).Value = -34988

$ws.Cells.Item(82, 8).Value = 1691.8334
$ws.Cells.Item(82, 9).Value = 1040
$ws.Cells.Item(82, 11).Value = 1040
$ws.Cells.Item(82, 13).Value = -679

$ws.Cells.Item(85, 8).Value = 1691.8334
$ws.Cells.Item(85, 9).Value = 1040
$ws.Cells.Item(85, 11).Value = 1040
$ws.Cells.Item(85, 13).Value = 208

$ws.Cells.Item(93, 8).Value = 475859.2
$ws.Cells.Item(93, 9).Value = 1651.7727
$ws.Cells.Item(93, 11).Value = 1651.7727
$ws.Cells.Item(93, 13).Value = -403.7727

$ws.Cells.Item(109, 8).Value = 65000
$ws.Cells.Item(109, 9).Value = 65000
$ws.Cells.Item(109, 11).Value = 65000
$ws.Cells.Item(109, 13).Value = -63613

$ws.Cells.Item(122, 8).Value = 8932857
$ws.Cells.Item(122, 9).Value = 3999.9
$ws.Cells.Item(122, 10).Value = 31255000
$ws.Cells.Item(122, 11).Value = 11999.7
$ws.Cells.Item(122, 12).Value = 93765000
$ws.Cells.Item(122, 13).Value = -9549.700000000001
$ws.Cells.Item(122, 14).Value = -93769900

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(50, 8).Value = 76481.336
$ws.Cells.Item(50, 10).Value = 76481.336
$ws.Cells.Item(50, 12).Value = 76481.336
$ws.Cells.Item(50, 14).Value = -77743.336

$ws.Cells.Item(51, 8).Value = 19500
$ws.Cells.Item(51, 9).Value = 19500
$ws.Cells.Item(51, 11).Value = 19500
$ws.Cells.Item(51, 13).Value = -18990

$ws.Cells.Item(52, 8).Value = 35247.5
$ws.Cells.Item(52, 9).Value = 30000
$ws.Cells.Item(52, 11).Value = 30000
$ws.Cells.Item(52, 13).Value = -29774

$ws.Cells.Item(107, 8).Value = 1370.3715
$ws.Cells.Item(107, 9).Value = 1087.05
$ws.Cells.Item(107, 11).Value = 3261.15
$ws.Cells.Item(107, 13).Value = -1341.15

$ws.Cells.Item(122, 8).Value = 22230034
$ws.Cells.Item(122, 9).Value = 4160.6
$ws.Cells.Item(122, 10).Value = 50012376
$ws.Cells.Item(122, 11).Value = 12481.8
$ws.Cells.Item(122, 12).Value = 150037128
$ws.Cells.Item(122, 13).Value = -10031.8
$ws.Cells.Item(122, 14).Value = -150042028

$ws.Cells.Item(125, 8).Value = 78333
$ws.Cells.Item(125, 10).Value = 78333
$ws.Cells.Item(125, 12).Value = 78333
$ws.Cells.Item(125, 14).Value = -88173
